# Update LDLC prices history
#
# A new price-check pass was run, so a new timestamped column (G) is
# appended after the existing history columns (D, E, F). The header row
# gets the new check's timestamp, and every product row gets the latest
# observed price copied into the new column (products with no tracked
# price yet stay blank, matching the existing D:F pattern for those rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastDataRow = 100
$lastRow = 204
$newCol = 7      # column G
$lastCol = 6     # column F (previous last "timestamp" column)

# --- Header (row 1): reuse F1's formatting for the new column, then set
#     its own value to the new check's timestamp.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(1, $newCol).Value2 = "2026-01-27 20:12:47"

# --- Data rows: carry the last known price (column F) into the new
#     column G for every product that currently has a tracked price.
for ($r = 2; $r -le $lastDataRow; $r++) {
    $ws.Cells.Item($r, $newCol).Value2 = $ws.Cells.Item($r, $lastCol).Value2
}

# --- Remaining rows (no price history yet): materialize a blank cell in
#     the new column too, same as the existing empty D/E/F cells there.
for ($r = ($lastDataRow + 1); $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $newCol).Style = "Normal"
}
